# Workbook/worksheet handles
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "RPI shock tube"
$ws2 = $wb.Worksheets.Item(2)   # "NUIG shock tube"

# ---------------------------------------------------------------------------
# NUIG shock tube (sheet2): the "P [atm]" and "T [K]" columns had been
# entered in the wrong order. Swap columns B and C - both the values AND
# their per-cell formatting - for the header row and all 30 data rows
# (rows 1-31) so the data lines up with the correct ChemKED headers.
# A scratch cell (Z1, outside the sheet's used range) is used to hold one
# side of the swap; Copy() carries formatting as well as the value.
# ---------------------------------------------------------------------------
$tmp = $ws2.Range("Z1")
for ($r = 1; $r -le 31; $r++) {
    $bCell = $ws2.Cells.Item($r, 2)
    $cCell = $ws2.Cells.Item($r, 3)
    [void]$bCell.Copy($tmp)
    [void]$cCell.Copy($bCell)
    [void]$tmp.Copy($cCell)
}
$tmp.Clear()

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state:
# NUIG shock tube selection moved to G8, RPI shock tube selection moved to
# E26 (and scrolled back to the top), with RPI shock tube left as the
# active (tab-selected) sheet.
# ---------------------------------------------------------------------------
$ws2.Activate()
[void]$ws2.Range("G8").Select()

$ws1.Activate()
[void]$ws1.Range("E26").Select()
